$wb = $excel.ActiveWorkbook

# NOTE: worksheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) lookups are case-insensitive, so sheets are addressed
# by their (1-based) tab position instead of by name to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# ---- Restricciones_del_follower (sheet 3): regenerated follower restriction table ----
$ws3 = $wb.Worksheets.Item(3)

# Columns B, D, E and F hold numbers-as-text; force Text formatting first so
# Excel does not auto-coerce the assigned strings back into numeric cells.
$ws3.Range("B2:B6").NumberFormat = "@"
$ws3.Range("D2:D6").NumberFormat = "@"
$ws3.Range("E2:E6").NumberFormat = "@"
$ws3.Range("F2:F6").NumberFormat = "@"

$ws3.Range("A2").Value = "-0.6941935483870978 - x + 2.5806451612903225y"
$ws3.Range("B2").Value = "2.694193548387098"
$ws3.Range("D2").Value = "0.21"
$ws3.Range("E2").Value = "3.8"
$ws3.Range("F2").Value = "2.4"

$ws3.Range("A3").Value = "-20.624086021505377 - 0.25x + 8.924731182795698y"
$ws3.Range("B3").Value = "18.624086021505377"
$ws3.Range("D3").Value = "0.19"
$ws3.Range("E3").Value = "9.5"
$ws3.Range("F3").Value = "8.299999999999999"

$ws3.Range("A4").Value = "-5.76645 + x + 0.03500000000000003y"
$ws3.Range("B4").Value = "-2.23355"
$ws3.Range("D4").Value = "0.97"
$ws3.Range("E4").Value = "0"
$ws3.Range("F4").Value = "5.0"

$ws3.Range("A5").Value = "-20.28193548387097 + x + 5.806451612903225y"
$ws3.Range("B5").Value = "18.02193548387097"
$ws3.Range("D5").Value = "0.43"
$ws3.Range("E5").Value = "1.5"
$ws3.Range("F5").Value = "5.4"

$ws3.Range("A6").Value = "-1.8991397849462368 + 0.7526881720430108y"
$ws3.Range("B6").Value = "1.8591397849462368"
$ws3.Range("D6").Value = "0.47"
$ws3.Range("E6").Value = "3.7"
$ws3.Range("F6").Value = "0.7000000000000001"

# ---- Punto_modificado (sheet 4) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2:B2").NumberFormat = "@"
$ws4.Range("A2").Value = "5.68"
$ws4.Range("B2").Value = "2.47"

# ---- Vector_bf (sheet 5) ----
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "-4.1221220430107515"

# ---- Vector_BF (sheet 6) ----
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2:A3").NumberFormat = "@"
$ws6.Range("A2").Value = "3.675"
$ws6.Range("A3").Value = "-107.08602150537634"

# ---- Vector_Alpha (sheet 7): A2 is a genuine numeric cell ----
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 0.9299999999999999
